# Auto-generated edit script applying the symbol-list refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving numeric-looking text (prices, %, hour) must be pre-
# formatted as Text so Excel stores the exact literal string instead of
# silently converting it to a number/percentage.
$numericTextCells = @(
    "D2",
    "E2",
    "G2",
    "D3",
    "E3",
    "G3",
    "D4",
    "E4",
    "G4",
    "D5",
    "E5",
    "G5",
    "D6",
    "E6",
    "G6",
    "D7",
    "E7",
    "G7",
    "D8",
    "E8",
    "G8",
    "D9",
    "E9",
    "G9",
    "D10",
    "E10",
    "G10",
    "D11",
    "E11",
    "G11",
    "D12",
    "E12",
    "G12",
    "D13",
    "E13",
    "G13",
    "D14",
    "E14",
    "G14",
    "D15",
    "E15",
    "G15",
    "D16",
    "E16",
    "G16",
    "D17",
    "E17",
    "G17",
    "D18",
    "E18",
    "G18",
    "E19",
    "G19",
    "D20",
    "E20",
    "G20",
    "D21",
    "E21",
    "G21",
    "D22",
    "E22",
    "G22",
    "D23",
    "E23",
    "G23",
    "D24",
    "E24",
    "G24",
    "D25",
    "E25",
    "G25",
    "D26",
    "E26",
    "G26",
    "D27",
    "E27",
    "G27",
    "G28",
    "G29",
    "G30",
    "G31",
    "G32",
    "G33",
    "G34",
    "G35",
    "G36",
    "G37",
    "G38",
    "G39",
    "D40",
    "E40",
    "G40",
    "D41",
    "E41",
    "G41",
    "D42",
    "E42",
    "G42",
    "D43",
    "E43",
    "G43",
    "D44",
    "E44",
    "G44",
    "D45",
    "E45",
    "G45",
    "D46",
    "E46",
    "G46",
    "D47",
    "E47",
    "G47",
    "D48",
    "E48",
    "G48",
    "D49",
    "E49",
    "G49",
    "D50",
    "E50",
    "G50",
    "G51"
)
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the new values (matches the commit's refreshed snapshot).
$ws.Range("D2").Value = '256.21'
$ws.Range("E2").Value = '0.50%'
$ws.Range("G2").Value = '15'
$ws.Range("D3").Value = '27.03'
$ws.Range("E3").Value = '-3.83%'
$ws.Range("G3").Value = '15'
$ws.Range("D4").Value = '4.585'
$ws.Range("E4").Value = '-11.93%'
$ws.Range("G4").Value = '15'
$ws.Range("D5").Value = '0.05896'
$ws.Range("E5").Value = '0.49%'
$ws.Range("G5").Value = '15'
$ws.Range("D6").Value = '6.644'
$ws.Range("E6").Value = '-0.62%'
$ws.Range("G6").Value = '15'
$ws.Range("D7").Value = '0.8675'
$ws.Range("E7").Value = '-0.21%'
$ws.Range("G7").Value = '15'
$ws.Range("D8").Value = '0.9388'
$ws.Range("E8").Value = '-1.96%'
$ws.Range("G8").Value = '15'
$ws.Range("D9").Value = '0.1401'
$ws.Range("E9").Value = '-0.60%'
$ws.Range("G9").Value = '15'
$ws.Range("D10").Value = '0.03831'
$ws.Range("E10").Value = '11.40%'
$ws.Range("G10").Value = '15'
$ws.Range("D11").Value = '0.07076'
$ws.Range("E11").Value = '-1.28%'
$ws.Range("G11").Value = '15'
$ws.Range("D12").Value = '0.03200'
$ws.Range("E12").Value = '-0.24%'
$ws.Range("G12").Value = '15'
$ws.Range("D13").Value = '0.09252'
$ws.Range("E13").Value = '0.38%'
$ws.Range("G13").Value = '15'
$ws.Range("D14").Value = '0.001545'
$ws.Range("E14").Value = '-0.22%'
$ws.Range("G14").Value = '15'
$ws.Range("D15").Value = '0.0006011'
$ws.Range("E15").Value = '-1.13%'
$ws.Range("G15").Value = '15'
$ws.Range("D16").Value = '0.006018'
$ws.Range("E16").Value = '3.56%'
$ws.Range("G16").Value = '15'
$ws.Range("D17").Value = '3.513'
$ws.Range("E17").Value = '0.48%'
$ws.Range("G17").Value = '15'
$ws.Range("D18").Value = '3.188'
$ws.Range("E18").Value = '-0.73%'
$ws.Range("G18").Value = '15'
$ws.Range("E19").Value = '-1.06%'
$ws.Range("G19").Value = '15'
$ws.Range("D20").Value = '0.3075'
$ws.Range("E20").Value = '-3.23%'
$ws.Range("G20").Value = '15'
$ws.Range("D21").Value = '0.1282'
$ws.Range("E21").Value = '-2.04%'
$ws.Range("G21").Value = '15'
$ws.Range("D22").Value = '3.842'
$ws.Range("E22").Value = '8.25%'
$ws.Range("G22").Value = '15'
$ws.Range("D23").Value = '0.04223'
$ws.Range("E23").Value = '1.20%'
$ws.Range("G23").Value = '15'
$ws.Range("B24").Value = 'BitKan'
$ws.Range("C24").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D24").Value = '0.001218'
$ws.Range("E24").Value = '-0.30%'
$ws.Range("G24").Value = '15'
$ws.Range("B25").Value = 'HotbitToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D25").Value = '0.004278'
$ws.Range("E25").Value = '-6.18%'
$ws.Range("G25").Value = '15'
$ws.Range("B26").Value = 'NitroEx'
$ws.Range("C26").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D26").Value = '0.0001198'
$ws.Range("E26").Value = '-0.16%'
$ws.Range("G26").Value = '15'
$ws.Range("B27").Value = 'UpBots'
$ws.Range("C27").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D27").Value = '0.0001505'
$ws.Range("E27").Value = '2.70%'
$ws.Range("G27").Value = '15'
$ws.Range("B28").Value = 'Spectre.aiUtilityToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut'
$ws.Range("D28").Value = '--'
$ws.Range("E28").Value = '--%'
$ws.Range("G28").Value = '15'
$ws.Range("B29").Value = 'LegolasExchange'
$ws.Range("C29").Value = 'https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo'
$ws.Range("G29").Value = '15'
$ws.Range("B30").Value = 'BitZToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz'
$ws.Range("G30").Value = '15'
$ws.Range("B31").Value = 'Birake'
$ws.Range("C31").Value = 'https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir'
$ws.Range("G31").Value = '15'
$ws.Range("B32").Value = 'ZBToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("G32").Value = '15'
$ws.Range("G33").Value = '15'
$ws.Range("G34").Value = '15'
$ws.Range("G35").Value = '15'
$ws.Range("G36").Value = '15'
$ws.Range("G37").Value = '15'
$ws.Range("G38").Value = '15'
$ws.Range("G39").Value = '15'
$ws.Range("D40").Value = '0.03815'
$ws.Range("E40").Value = '-0.13%'
$ws.Range("G40").Value = '15'
$ws.Range("D41").Value = '0.006212'
$ws.Range("E41").Value = '61.36%'
$ws.Range("G41").Value = '15'
$ws.Range("D42").Value = '0.1098'
$ws.Range("E42").Value = '-0.39%'
$ws.Range("G42").Value = '15'
$ws.Range("D43").Value = '0.002283'
$ws.Range("E43").Value = '-7.20%'
$ws.Range("G43").Value = '15'
$ws.Range("D44").Value = '0.01154'
$ws.Range("E44").Value = '18.58%'
$ws.Range("G44").Value = '15'
$ws.Range("D45").Value = '0.00005463'
$ws.Range("E45").Value = '1.45%'
$ws.Range("G45").Value = '15'
$ws.Range("D46").Value = '0.00000000749'
$ws.Range("E46").Value = '-0.12%'
$ws.Range("G46").Value = '15'
$ws.Range("D47").Value = '0.07769'
$ws.Range("E47").Value = '-13.66%'
$ws.Range("G47").Value = '15'
$ws.Range("D48").Value = '0.002275'
$ws.Range("E48").Value = '6.87%'
$ws.Range("G48").Value = '15'
$ws.Range("D49").Value = '0.00002097'
$ws.Range("E49").Value = '-0.12%'
$ws.Range("G49").Value = '15'
$ws.Range("D50").Value = '0.0001997'
$ws.Range("E50").Value = '-0.12%'
$ws.Range("G50").Value = '15'
$ws.Range("G51").Value = '15'
